$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 60.333332
$ws.Range("I4").Value = 80.5
$ws.Range("K4").Value = 80.5
$ws.Range("M4").Value = 33.5
$ws.Range("H17").Value = 2836.4912
$ws.Range("J17").Value = 2793.3928
$ws.Range("L17").Value = 8380.178400000001
$ws.Range("N17").Value = -8716.178400000001
$ws.Range("H19").Value = 1218.5264
$ws.Range("I19").Value = 1237.9333
$ws.Range("J19").Value = 1145.75
$ws.Range("K19").Value = 1237.9333
$ws.Range("L19").Value = 1145.75
$ws.Range("M19").Value = -1062.9333
$ws.Range("N19").Value = -1495.75
$ws.Range("H64").Value = 5824.8
$ws.Range("I64").Value = 4748
$ws.Range("K64").Value = 4748
$ws.Range("M64").Value = -4500
$ws.Range("H67").Value = 5824.8
$ws.Range("I67").Value = 4748
$ws.Range("K67").Value = 4748
$ws.Range("M67").Value = -3890
$ws.Range("H132").Value = 2572.3684
$ws.Range("I132").Value = 2448.611
$ws.Range("K132").Value = 7345.833
$ws.Range("M132").Value = -4815.833
$ws.Range("H137").Value = 3473538.8
$ws.Range("I137").Value = 1011.5
$ws.Range("J137").Value = 5557055
$ws.Range("K137").Value = 3034.5
$ws.Range("L137").Value = 16671165
$ws.Range("M137").Value = -484.5
$ws.Range("N137").Value = -16676265

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5265.5557
$ws.Range("I45").Value = 4722.5
$ws.Range("J45").Value = 5700
$ws.Range("K45").Value = 4722.5
$ws.Range("L45").Value = 5700
$ws.Range("M45").Value = -4345.5
$ws.Range("N45").Value = -6454
$ws.Range("H62").Value = 9000000
$ws.Range("J62").Value = 9000000
$ws.Range("L62").Value = 9000000
$ws.Range("N62").Value = -9001248
$ws.Range("H65").Value = 9000000
$ws.Range("J65").Value = 9000000
$ws.Range("L65").Value = 27000000
$ws.Range("N65").Value = -27006240
$ws.Range("H122").Value = 3112.6924
$ws.Range("I122").Value = 2916.3333
$ws.Range("J122").Value = 3380.4546
$ws.Range("K122").Value = 8748.999899999999
$ws.Range("L122").Value = 10141.3638
$ws.Range("M122").Value = -6298.999899999999
$ws.Range("N122").Value = -15041.3638
$ws.Range("H132").Value = 2750.6897
$ws.Range("I132").Value = 2664.6
$ws.Range("J132").Value = 3288.75
$ws.Range("K132").Value = 7993.799999999999
$ws.Range("L132").Value = 9866.25
$ws.Range("M132").Value = -5463.799999999999
$ws.Range("N132").Value = -14926.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 48988.184
$ws.Range("I20").Value = 85980.5
$ws.Range("J20").Value = 4597.4
$ws.Range("K20").Value = 85980.5
$ws.Range("L20").Value = 4597.4
$ws.Range("M20").Value = -85733.5
$ws.Range("N20").Value = -5091.4
$ws.Range("H130").Value = 178000
$ws.Range("J130").Value = 178000
$ws.Range("L130").Value = 178000
$ws.Range("N130").Value = -188040
$ws.Range("H134").Value = 3091.45
$ws.Range("I134").Value = 2238.4375
$ws.Range("J134").Value = 6503.5
$ws.Range("K134").Value = 6715.3125
$ws.Range("L134").Value = 19510.5
$ws.Range("M134").Value = -4180.3125
$ws.Range("N134").Value = -24580.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3352
$ws.Range("I31").Value = 813.93335
$ws.Range("K31").Value = 813.93335
$ws.Range("M31").Value = -518.93335
$ws.Range("H34").Value = 3352
$ws.Range("I34").Value = 813.93335
$ws.Range("K34").Value = 813.93335
$ws.Range("M34").Value = -611.93335
$ws.Range("H36").Value = 47009.4
$ws.Range("I36").Value = 47009.4
$ws.Range("K36").Value = 47009.4
$ws.Range("M36").Value = -46621.4
$ws.Range("H40").Value = 47009.4
$ws.Range("I40").Value = 47009.4
$ws.Range("K40").Value = 47009.4
$ws.Range("M40").Value = -46849.4
$ws.Range("H58").Value = 2807.7334
$ws.Range("J58").Value = 3721.1428
$ws.Range("L58").Value = 3721.1428
$ws.Range("N58").Value = -4127.1428
$ws.Range("H62").Value = 3179.6
$ws.Range("J62").Value = 2966
$ws.Range("L62").Value = 2966
$ws.Range("N62").Value = -4214
$ws.Range("H65").Value = 3179.6
$ws.Range("J65").Value = 2966
$ws.Range("L65").Value = 14830
$ws.Range("N65").Value = -21070
$ws.Range("H86").Value = 33785.2
$ws.Range("I86").Value = 35997.25
$ws.Range("K86").Value = 35997.25
$ws.Range("M86").Value = -34874.25
$ws.Range("H87").Value = 96994.5
$ws.Range("J87").Value = 96994.5
$ws.Range("L87").Value = 96994.5
$ws.Range("N87").Value = -99366.5
$ws.Range("H88").Value = 45196.668
$ws.Range("J88").Value = 45196.668
$ws.Range("L88").Value = 45196.668
$ws.Range("N88").Value = -46008.668
$ws.Range("H89").Value = 33785.2
$ws.Range("I89").Value = 35997.25
$ws.Range("K89").Value = 179986.25
$ws.Range("M89").Value = -174370.25
$ws.Range("H90").Value = 96994.5
$ws.Range("J90").Value = 96994.5
$ws.Range("L90").Value = 290983.5
$ws.Range("N90").Value = -302839.5
$ws.Range("H91").Value = 45196.668
$ws.Range("J91").Value = 45196.668
$ws.Range("L91").Value = 45196.668
$ws.Range("N91").Value = -48004.668
$ws.Range("H134").Value = 3200
$ws.Range("I134").Value = 2933.3333
$ws.Range("K134").Value = 8799.999899999999
$ws.Range("M134").Value = -6264.999899999999
$ws.Range("H136").Value = 2807.7334
$ws.Range("J136").Value = 3721.1428
$ws.Range("L136").Value = 11163.4284
$ws.Range("N136").Value = -16263.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 12554.625
$ws.Range("I2").Value = 30.818182
$ws.Range("J2").Value = 40107
$ws.Range("K2").Value = 30.818182
$ws.Range("L2").Value = 40107
$ws.Range("M2").Value = 82.18181799999999
$ws.Range("N2").Value = -40333
$ws.Range("H70").Value = 22942.06
$ws.Range("I70").Value = 91469.71000000001
$ws.Range("J70").Value = 4492.3076
$ws.Range("K70").Value = 91469.71000000001
$ws.Range("L70").Value = 4492.3076
$ws.Range("M70").Value = -91199.71000000001
$ws.Range("N70").Value = -5032.3076
$ws.Range("H73").Value = 22942.06
$ws.Range("I73").Value = 91469.71000000001
$ws.Range("J73").Value = 4492.3076
$ws.Range("K73").Value = 91469.71000000001
$ws.Range("L73").Value = 4492.3076
$ws.Range("M73").Value = -90533.71000000001
$ws.Range("N73").Value = -6364.3076
$ws.Range("H132").Value = 2899.4546
$ws.Range("I132").Value = 2639.4
$ws.Range("K132").Value = 7918.200000000001
$ws.Range("M132").Value = -5388.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2039.5333
$ws.Range("I16").Value = 1476.6923
$ws.Range("J16").Value = 5698
$ws.Range("K16").Value = 1476.6923
$ws.Range("L16").Value = 5698
$ws.Range("M16").Value = -1306.6923
$ws.Range("N16").Value = -6038
$ws.Range("H22").Value = 3916.889
$ws.Range("I22").Value = 3100.6
$ws.Range("J22").Value = 4937.25
$ws.Range("K22").Value = 3100.6
$ws.Range("L22").Value = 4937.25
$ws.Range("M22").Value = -2805.6
$ws.Range("N22").Value = -5527.25
$ws.Range("H27").Value = 3916.889
$ws.Range("I27").Value = 3100.6
$ws.Range("J27").Value = 4937.25
$ws.Range("K27").Value = 3100.6
$ws.Range("L27").Value = 4937.25
$ws.Range("M27").Value = -2993.6
$ws.Range("N27").Value = -5151.25
$ws.Range("H46").Value = 9507.105
$ws.Range("I46").Value = 4000.5
$ws.Range("J46").Value = 10154.941
$ws.Range("K46").Value = 4000.5
$ws.Range("L46").Value = 10154.941
$ws.Range("M46").Value = -3812.5
$ws.Range("N46").Value = -10530.941
$ws.Range("H93").Value = 3006.125
$ws.Range("I93").Value = 2300.5
$ws.Range("K93").Value = 2300.5
$ws.Range("M93").Value = -1052.5
$ws.Range("H136").Value = 20000.572
$ws.Range("I136").Value = 10002
$ws.Range("J136").Value = 24000
$ws.Range("K136").Value = 30006
$ws.Range("L136").Value = 72000
$ws.Range("M136").Value = -27456
$ws.Range("N136").Value = -77100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3249.1
$ws.Range("I132").Value = 3199
$ws.Range("K132").Value = 9597
$ws.Range("M132").Value = -7067
